$d = $word.ActiveDocument

# --- Block 1: "Nếu chạm Layer chặn" -> three ilvl=1 bullets ---------------
$d.Paragraphs(2).Range.Text = "Gọi điểm tại vị trí bắt đầu bị chặn là a"

$null = $d.Paragraphs(2).Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.Text = "Gọi điểm đầu tiên thoát khỏi bị chặn là b"

$null = $d.Paragraphs(3).Range.InsertParagraphAfter()
$d.Paragraphs(4).Range.Text = "Sử dụng điểm ở giữa a và b"

# --- Block 2: "Nếu vẽ từ trên xuống/từ dưới lên" -> "Có bug" + children --
$d.Paragraphs(5).Range.Text = "Có bug"

# drop its two old ilvl=3 children
$d.Paragraphs(6).Range.Delete()
$d.Paragraphs(6).Range.Delete()

# add three new ilvl=1 bullets after "Có bug"
$null = $d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(6).Range.Text = "Cần xét tất cả các điểm trên đường thẳng từ a đến b"

$null = $d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = "Nhấn mạnh là vị trí thoát đầu tiên nhá, vì điểm vị chặn đầu tiên luôn nằm sát vớ"
$d.Paragraphs(7).Range.InsertAfter("i Layer block")

$null = $d.Paragraphs(7).Range.InsertParagraphAfter()
$d.Paragraphs(8).Range.Text = "Thử thêm collider trước nếu va chạm thì dừng"

# --- Block 3: "Nếu vẽ từ trái/phải sang" -> "Tạo một đối tượng rỗng" -----
$d.Paragraphs(9).Range.Text = "Tạo một đối tượng rỗng"

# drop its two old ilvl=3 children
$d.Paragraphs(10).Range.Delete()
$d.Paragraphs(10).Range.Delete()

# add one new ilvl=1 bullet after it
$null = $d.Paragraphs(9).Range.InsertParagraphAfter()
$d.Paragraphs(10).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(10).Range.Text = "Hoặc tạo 1 Line từ a đến b"

# --- Block 4: "Làm sao xác định đc vẽ từ hướng nào??" --------------------
$d.Paragraphs(11).Range.Text = "Nếu va chạm với Blocked thì Line chính k vẽ"

# drop its two old ilvl=3 children ("Gọi a=...", "Gọi b=...")
$d.Paragraphs(12).Range.Delete()
$d.Paragraphs(12).Range.Delete()

# --- Block 5: "Lấy điểm ở giữa a và b xe có bị chặn k" --------------------
# Drop the two plain paragraphs after it ("Lấy mousePosition liên tục" and
# "Nếu thảo mãn các yếu tổ ở trên thì nối Line lại") -- neither carries the
# bookmark, so a normal delete is safe.
$d.Paragraphs(13).Range.Delete()
$d.Paragraphs(13).Range.Delete()

# Paragraph 13 now is "Vẽ ngang ok, vẽ dọc có bug" -- the _GoBack bookmark
# sits right after its run, inside the same <w:p>. Remove only its text
# (not the whole paragraph/paragraph-mark) so the bookmark survives, then
# fold what remains into paragraph 12 by deleting paragraph 12's own
# paragraph mark -- this merges the (now empty) bookmark-bearing paragraph
# into paragraph 12, landing the bookmark right after paragraph 12's text.
$d.Paragraphs(12).Range.Text = "Ngược lại thì Line chính vẽ"

$p13 = $d.Paragraphs(13)
$r13 = $p13.Range
$textOnly = $d.Range($r13.Start, $r13.End - 1)
$textOnly.Delete()

$p12 = $d.Paragraphs(12)
$r12 = $p12.Range
$mark = $d.Range($r12.End - 1, $r12.End)
$mark.Delete()

# --- Block 6: add two new bullets at the tail -----------------------------
# Paragraphs 13 = "Thiết kế xe", 14 = "Khi nhả nút...", 15 = "Nhớ thêm
# CantDrawOver..." are all unchanged; append after the last of these.
$null = $d.Paragraphs(15).Range.InsertParagraphAfter()
$d.Paragraphs(16).Range.Text = "Tối ưu giảm xóc"

$null = $d.Paragraphs(16).Range.InsertParagraphAfter()
$d.Paragraphs(17).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(17).Range.Text = "Có thể thử kích hoạt sử dụng collider"
